$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.098.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.922.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.62%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.75%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5254"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4059"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.131"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.424"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.919.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.409"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("E16").Value = "  -0.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.19%  "

$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06703"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.062"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.083.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.140.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.80%  "

$ws.Range("E27").Value = "  +2.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.469"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.082"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1059"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.115"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.651"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02527"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06600"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("E37").Value = "  +3.44%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.118"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.49%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.242"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.213"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6573"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.248"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6211"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.790"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.092"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.250"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.95%  "

$ws.Range("B50").Value = "WEMIXTOKEN"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.160"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.34%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.93%  "
